# Refresh Universalis market-price snapshots (columns H:N) on every
# Leve-profit worksheet. Values below are the latest pull from the
# scheduled pricing runner; row/column layout is unchanged.
$wb = $excel.ActiveWorkbook

# --- ALC sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
# row 33
$ws.Range("H33").Value = 245.05556
$ws.Range("I33").Value = 245.05556
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 245.05556
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -16.05556000000001
$ws.Range("N33").ClearContents()
# row 53
$ws.Range("H53").Value = 1149.1765
$ws.Range("J53").Value = 1383.1428
$ws.Range("L53").Value = 1383.1428
$ws.Range("N53").Value = -2657.1428
# row 61
$ws.Range("H61").Value = 4759.421
$ws.Range("I61").Value = 4759.421
$ws.Range("K61").Value = 14278.263
$ws.Range("M61").Value = -14106.263
# row 100
$ws.Range("H100").Value = 1952
$ws.Range("I100").Value = 2057.7778
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 2057.7778
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -1516.7778
$ws.Range("N100").Value = -2082
# row 132
$ws.Range("H132").Value = 2289.681
$ws.Range("J132").Value = 2884.8572
$ws.Range("L132").Value = 8654.571599999999
$ws.Range("N132").Value = -13714.5716
# row 138
$ws.Range("H138").Value = 2466.8438
$ws.Range("J138").Value = 2668.087
$ws.Range("L138").Value = 8004.261
$ws.Range("N138").Value = -18284.261

# --- ARM sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 6158.094
$ws.Range("I32").Value = 3240.0945
$ws.Range("K32").Value = 3240.0945
$ws.Range("M32").Value = -2953.0945
# row 63
$ws.Range("H63").Value = 2127.5
$ws.Range("J63").Value = 349
$ws.Range("L63").Value = 349
$ws.Range("N63").Value = -1721
# row 66
$ws.Range("H66").Value = 2127.5
$ws.Range("J66").Value = 349
$ws.Range("L66").Value = 1745
$ws.Range("N66").Value = -8609

# --- BSM sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Range("H86").Value = 2347.2222
$ws.Range("I86").Value = 2047.579
$ws.Range("J86").Value = 3058.875
$ws.Range("K86").Value = 2047.579
$ws.Range("L86").Value = 3058.875
$ws.Range("M86").Value = -924.579
$ws.Range("N86").Value = -5304.875
# row 89
$ws.Range("H89").Value = 2347.2222
$ws.Range("I89").Value = 2047.579
$ws.Range("J89").Value = 3058.875
$ws.Range("K89").Value = 10237.895
$ws.Range("L89").Value = 15294.375
$ws.Range("M89").Value = -4621.895
$ws.Range("N89").Value = -26526.375

# --- CRP sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 266443.03
$ws.Range("I31").Value = 360514.22
$ws.Range("K31").Value = 360514.22
$ws.Range("M31").Value = -360219.22
# row 34
$ws.Range("H34").Value = 266443.03
$ws.Range("I34").Value = 360514.22
$ws.Range("K34").Value = 360514.22
$ws.Range("M34").Value = -360312.22
# row 135
$ws.Range("H135").Value = 86666
$ws.Range("I135").Value = 60000
$ws.Range("K135").Value = 60000
$ws.Range("M135").Value = -54930

# --- CUL sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
# row 22
$ws.Range("H22").Value = 1000
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 3000
$ws.Range("N22").Value = -3338
# row 27
$ws.Range("H27").Value = 1000
$ws.Range("J27").Value = 1000
$ws.Range("L27").Value = 3000
$ws.Range("N27").Value = -3204
# row 86
$ws.Range("H86").Value = 444.0909
$ws.Range("I86").Value = 440.94446
$ws.Range("J86").Value = 458.25
$ws.Range("K86").Value = 1322.83338
$ws.Range("L86").Value = 1374.75
$ws.Range("M86").Value = -136.83338
$ws.Range("N86").Value = -3746.75
# row 89
$ws.Range("H89").Value = 444.0909
$ws.Range("I89").Value = 440.94446
$ws.Range("J89").Value = 458.25
$ws.Range("K89").Value = 3968.50014
$ws.Range("L89").Value = 4124.25
$ws.Range("M89").Value = 1959.49986
$ws.Range("N89").Value = -15980.25
# row 108
$ws.Range("H108").Value = 1171.8334
$ws.Range("I108").Value = 1171.8334
$ws.Range("K108").Value = 3515.5002
$ws.Range("M108").Value = -635.5001999999999

# --- GSM sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
# row 80
$ws.Range("H80").Value = 3862.919
$ws.Range("I80").Value = 3585.8125
$ws.Range("J80").Value = 4074.0476
$ws.Range("K80").Value = 3585.8125
$ws.Range("L80").Value = 4074.0476
$ws.Range("M80").Value = -2587.8125
$ws.Range("N80").Value = -6070.0476
# row 83
$ws.Range("H83").Value = 3862.919
$ws.Range("I83").Value = 3585.8125
$ws.Range("J83").Value = 4074.0476
$ws.Range("K83").Value = 17929.0625
$ws.Range("L83").Value = 20370.238
$ws.Range("M83").Value = -12937.0625
$ws.Range("N83").Value = -30354.238

# --- LTW sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 4298.1
$ws.Range("I7").Value = 3948.3333
$ws.Range("J7").Value = 4822.75
$ws.Range("K7").Value = 3948.3333
$ws.Range("L7").Value = 4822.75
$ws.Range("M7").Value = -3836.3333
$ws.Range("N7").Value = -5046.75
# row 40
$ws.Range("H40").Value = 5615.385
$ws.Range("I40").Value = 4817.727
$ws.Range("K40").Value = 4817.727
$ws.Range("M40").Value = -4681.727
# row 126
$ws.Range("H126").Value = 4298.1
$ws.Range("I126").Value = 3948.3333
$ws.Range("J126").Value = 4822.75
$ws.Range("K126").Value = 11844.9999
$ws.Range("L126").Value = 14468.25
$ws.Range("M126").Value = -9374.999899999999
$ws.Range("N126").Value = -19408.25

# --- WVR sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 9182.462
$ws.Range("I62").Value = 7000
$ws.Range("J62").Value = 9579.272000000001
$ws.Range("K62").Value = 7000
$ws.Range("L62").Value = 9579.272000000001
$ws.Range("M62").Value = -6376
$ws.Range("N62").Value = -10827.272
# row 65
$ws.Range("H65").Value = 9182.462
$ws.Range("I65").Value = 7000
$ws.Range("J65").Value = 9579.272000000001
$ws.Range("K65").Value = 35000
$ws.Range("L65").Value = 47896.36
$ws.Range("M65").Value = -31880
$ws.Range("N65").Value = -54136.36
# row 81
$ws.Range("H81").Value = 5153.4546
$ws.Range("I81").Value = 2868.9
$ws.Range("J81").Value = 27999
$ws.Range("K81").Value = 5737.8
$ws.Range("L81").Value = 55998
$ws.Range("M81").Value = -4676.8
$ws.Range("N81").Value = -58120
# row 84
$ws.Range("H84").Value = 5153.4546
$ws.Range("I84").Value = 2868.9
$ws.Range("J84").Value = 27999
$ws.Range("K84").Value = 28689
$ws.Range("L84").Value = 279990
$ws.Range("M84").Value = -23385
$ws.Range("N84").Value = -290598
# row 96
$ws.Range("H96").Value = 39472.965
$ws.Range("I96").Value = 65362.438
$ws.Range("J96").Value = 1815.5454
$ws.Range("K96").Value = 65362.438
$ws.Range("L96").Value = 1815.5454
$ws.Range("M96").Value = -63989.438
$ws.Range("N96").Value = -4561.5454
# row 122
$ws.Range("H122").Value = 2549.7273
$ws.Range("I122").Value = 2550.3333
$ws.Range("J122").Value = 2547
$ws.Range("K122").Value = 7650.999899999999
$ws.Range("L122").Value = 7641
$ws.Range("M122").Value = -5200.999899999999
$ws.Range("N122").Value = -12541
# row 136
$ws.Range("H136").Value = 772281.4399999999
$ws.Range("J136").Value = 9000
$ws.Range("L136").Value = 27000
$ws.Range("N136").Value = -32100
